# "updated single child outputs" - refresh the two predicted scores on the
# quadratic-svm-score sheet and tidy up formatting that came along with the
# re-export (narrower duplicate text style collapses onto the shared one,
# column B widened slightly to fit the new decimal values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New prediction values for the two genome rows.
$ws.Range("B2").Value = 0.087143138009892596
$ws.Range("B3").Value = -0.04862773867910164

# Column B grew a touch to accommodate the longer decimal values.
$ws.Columns.Item(2).ColumnWidth = 13.7109375

# The header row and the "Row" label column were carrying a duplicate text
# style (same text format, blank border); re-applying the text format makes
# them collapse back onto the original shared style.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A3").NumberFormat = "@"
